$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "...through an interactive globe or map, making..." ->
#    "...through an interactive map and categorical data filters, making..."
# -----------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute(
    "through an interactive globe or map, making",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "through an interactive map and categorical data filters, making",
    2
)

# -----------------------------------------------------------------
# 2) Insert a brand-new paragraph right after the "Additionally, the
#    visualization will include current trade data..." paragraph,
#    describing the target user groups.
# -----------------------------------------------------------------
$additionallyPara = $d.Paragraphs.Item(5)
$insertionPoint = $additionallyPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(6)
$newParaStart = $newPara.Range.Start

# -----------------------------------------------------------------
# 3) Fill the new paragraph with its text (leading tab + two sentences
#    worth of content, matching the authored copy).
# -----------------------------------------------------------------
$newContent = $newPara.Range.Duplicate
$newContent.Collapse(1)
$newContent.InsertAfter("`tWe will be considering the wants and use cases of several user groups when designing our visualization, and considering how our data will be displayed. The first group we will consider is those that will be using our visualization to gather data and produce inferences or have some other practical use for the visualization. These users may include International business professionals, political scientists, and those studying international business. By focusing on the clear and effective presentation of our data, we will be able to make our visualization particularly useful to this user group.  Secondly, we hope to target also users that have a particular interest in the data we are using, but lack a practical purpose for it. Secondarily focusing on creating an elegant, intuitive interface will allow us to also focus on the needs of these users, while keeping our focus on our primary target group. By considering the specific needs of these users in the design phase, we hope to create a visualization that will also appeal to users that do not fall into our target groups, and who will be using our visualization for their own enjoyment. ")

# -----------------------------------------------------------------
# 4) Relocate the "_GoBack" bookmark from the end of the "Additionally…"
#    paragraph to the very start of the new paragraph we just wrote.
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bookmarkRange = $d.Range($newParaStart, $newParaStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
